# Update the "Price" (D) and "Volume(1h)" (E) columns for the latest crypto snapshot.
# Note: several D-column prices are numeric-looking text (e.g. "1.00", "1.40") that must
# stay stored as literal text (matching the source data's inline-string cells) instead of
# being auto-converted to numbers by Excel. A leading apostrophe (doubled here because the
# PowerShell string itself is single-quoted) forces those particular assignments to remain
# text while leaving every other cell/style in the workbook untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.243.07'
$ws.Range("E2").Value = '  -3.42%  '
$ws.Range("D3").Value = '2.368.43'
$ws.Range("E3").Value = '  -4.33%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '''499.48'
$ws.Range("E5").Value = '  -2.24%  '
$ws.Range("D6").Value = '''129.84'
$ws.Range("E6").Value = '  -3.02%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''0.542'
$ws.Range("E8").Value = '  -3.46%  '
$ws.Range("D9").Value = '2.373.50'
$ws.Range("E9").Value = '  -4.12%  '
$ws.Range("D10").Value = '''0.0981'
$ws.Range("E10").Value = '  -0.89%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").Value = '''0.324'
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("D14").Value = '2.788.83'
$ws.Range("E14").Value = '  -4.26%  '
$ws.Range("D15").Value = '56.194.62'
$ws.Range("E15").Value = '  -3.57%  '
$ws.Range("E16").Value = '  -2.99%  '
$ws.Range("E17").Value = '  -2.92%  '
$ws.Range("D18").Value = '2.306.59'
$ws.Range("E18").Value = '  -6.10%  '
$ws.Range("E19").Value = '  -4.08%  '
$ws.Range("E20").Value = '  -4.00%  '
$ws.Range("D21").Value = '''306.35'
$ws.Range("E21").Value = '  -3.84%  '
$ws.Range("D22").Value = '''6.23'
$ws.Range("E22").Value = '  -3.77%  '
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").Value = '''65.09'
$ws.Range("E24").Value = '  -1.32%  '
$ws.Range("D25").Value = '''0.996'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("D26").Value = '''0.367'
$ws.Range("E26").Value = '  -5.83%  '
$ws.Range("D27").Value = '''0.147'
$ws.Range("E27").Value = '  -5.26%  '
$ws.Range("E28").Value = '  -6.25%  '
$ws.Range("D29").Value = '''172.16'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("E30").Value = '  -4.47%  '
$ws.Range("E31").Value = '  -4.50%  '
$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").Value = '''0.998'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").Value = '''5.72'
$ws.Range("E34").Value = '  -8.97%  '
$ws.Range("E35").Value = '  -7.00%  '
$ws.Range("D36").Value = '''17.55'
$ws.Range("E36").Value = '  -3.28%  '
$ws.Range("E37").Value = '  -8.07%  '
$ws.Range("E38").Value = '  -3.96%  '
$ws.Range("D39").Value = '''36.01'
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("D40").Value = '''0.788'
$ws.Range("E40").Value = '  -4.82%  '
$ws.Range("D41").Value = '''1.40'
$ws.Range("E41").Value = '  -5.36%  '
$ws.Range("D42").Value = '''130.34'
$ws.Range("E42").Value = '  -5.63%  '
$ws.Range("E43").Value = '  -3.15%  '
$ws.Range("E44").Value = '  -4.85%  '
$ws.Range("D45").Value = '''0.565'
$ws.Range("E45").Value = '  -2.30%  '
$ws.Range("D46").Value = '''0.0903'
$ws.Range("E46").Value = '  -2.17%  '
$ws.Range("D47").Value = '''239.89'
$ws.Range("E47").Value = '  -9.19%  '
$ws.Range("D48").Value = '''0.0479'
$ws.Range("E48").Value = '  -4.77%  '
$ws.Range("E49").Value = '  -4.37%  '
$ws.Range("D50").Value = '''16.93'
$ws.Range("E50").Value = '  -3.49%  '
$ws.Range("E51").Value = '  -4.11%  '
